$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Describe Cuarentena" column (O10:O307) holds one shared string that is
# repeated for every municipality row. Updating the whole range to the new
# wording re-shares a single string entry (matching the original file's
# shared-string layout) instead of forking a brand new entry for a single
# cell.
$ws.Range("O10:O307").Value = "Segmentando a la poblacion para poder circular conforme a la terminacion de los digitos de su tarjeta de identidad, pasaporte o carnet de residente para extranjeros, para que puedan abastecerse de insumos básicos,  con horario de 5:00 am a 9:00 pm.  De lunes a domingo."

# Move the view: scroll so column F / row 4 is the top-left visible cell, and
# select K12 (as captured on the author's screen at save time).
$excel.ActiveWindow.ScrollColumn = 6
$excel.ActiveWindow.ScrollRow = 4
$ws.Range("K12").Select()
